$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (title changes from "Through 2022-07-17" to "Through 2022-07-18")
$ws.Name = "Through 2022-07-18"

# Update the row-label text for July in column A (through 07-17 -> through 07-18)
$ws.Range("A8").Value = "July (through 07-18)"

# Update July row (row 8) values for columns C..I (B unchanged)
$ws.Range("C8").Value = 37
$ws.Range("D8").Value = 37
$ws.Range("E8").Value = 43
$ws.Range("F8").Value = 27
$ws.Range("G8").Value = 74
$ws.Range("H8").Value = 89
$ws.Range("I8").Value = 99

# Update Total row (row 9) values for columns C..I (B unchanged)
$ws.Range("C9").Value = 285
$ws.Range("D9").Value = 427
$ws.Range("E9").Value = 396
$ws.Range("F9").Value = 278
$ws.Range("G9").Value = 546
$ws.Range("H9").Value = 849
$ws.Range("I9").Value = 904
